$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.170.09'
$ws.Range("E2").Value = '  +1.37%  '
$ws.Range("D3").Value = '1.613.95'
$ws.Range("E3").Value = '  +1.05%  '
$ws.Range("E4").Value = '  -0.52%  '
$ws.Range("D5").Value = '213.31'
$ws.Range("E5").Value = '  +2.30%  '
$ws.Range("E6").Value = '  -0.50%  '
$ws.Range("D7").Value = '0.483'
$ws.Range("E7").Value = '  +1.06%  '
$ws.Range("E8").Value = '  +1.82%  '
$ws.Range("E9").Value = '  +1.76%  '
$ws.Range("D10").Value = '18.42'
$ws.Range("E10").Value = '  +3.29%  '
$ws.Range("D11").Value = '0.0798'
$ws.Range("E11").Value = '  +1.32%  '
$ws.Range("D12").Value = '1.838.42'
$ws.Range("E12").Value = '  +0.99%  '
$ws.Range("D13").Value = '1.609.66'
$ws.Range("E13").Value = '  +0.71%  '
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("E15").Value = '  +0.93%  '
$ws.Range("D16").Value = '26.164.68'
$ws.Range("E16").Value = '  +1.32%  '
$ws.Range("D17").Value = '60.90'
$ws.Range("E17").Value = '  +0.89%  '
$ws.Range("E18").Value = '  +2.51%  '
$ws.Range("E19").Value = '  -0.39%  '
$ws.Range("D20").Value = '198.63'
$ws.Range("E20").Value = '  +5.04%  '
$ws.Range("E21").Value = '  +2.43%  '
$ws.Range("E22").Value = '  +2.22%  '
$ws.Range("E23").Value = '  +1.80%  '
$ws.Range("E24").Value = '  +2.67%  '
$ws.Range("D25").Value = '142.46'
$ws.Range("E25").Value = '  +0.80%  '
$ws.Range("E26").Value = '  +1.35%  '
$ws.Range("E27").Value = '  -0.51%  '
$ws.Range("E28").Value = '  +2.29%  '
$ws.Range("D29").Value = '6.52'
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("D30").Value = '1.18'
$ws.Range("E30").Value = '  -1.71%  '
$ws.Range("E31").Value = '  +3.41%  '
$ws.Range("D32").Value = '3.16'
$ws.Range("E32").Value = '  +2.71%  '
$ws.Range("E34").Value = '  +4.01%  '
$ws.Range("D35").Value = '2.35'
$ws.Range("E35").Value = '  -1.96%  '
$ws.Range("D36").Value = '1.108.00'
$ws.Range("E36").Value = '  +0.97%  '
$ws.Range("E37").Value = '  +1.80%  '
$ws.Range("E38").Value = '  -0.52%  '
$ws.Range("E39").Value = '  +2.78%  '
$ws.Range("D40").Value = '2.33'
$ws.Range("E40").Value = '  -1.15%  '
$ws.Range("D41").Value = '0.792'
$ws.Range("D42").Value = '0.799'
$ws.Range("E42").Value = '  +7.99%  '
$ws.Range("D43").Value = '1.750.50'
$ws.Range("E43").Value = '  +1.02%  '
$ws.Range("E44").Value = '  +1.25%  '
$ws.Range("D45").Value = '93.19'
$ws.Range("E45").Value = '  -2.48%  '
$ws.Range("D46").Value = '0.0₆0106'
$ws.Range("E46").Value = '  +5.55%  '
$ws.Range("E47").Value = '  +8.75%  '
$ws.Range("D48").Value = '54.05'
$ws.Range("E48").Value = '  +2.09%  '
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("D50").Value = '0.409'
$ws.Range("E50").Value = '  -0.30%  '
$ws.Range("E51").Value = '  -0.19%  '
